$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - admin user record
$ws.Range("A1").Value = "admin"
$ws.Range("B1").Value = "admin"
$ws.Range("C1").Value = "default"
$ws.Range("D1").Value = "'false"
$ws.Range("E1").Value = "admin"

# Row 2 - new user record
$ws.Range("A2").Value = "nuevo"
$ws.Range("B2").Value = "k"
$ws.Range("C2").Value = "admin-nuevo"
$ws.Range("D2").Value = "'false"
$ws.Range("E2").Value = "j"

# The leading apostrophe above forces the "false" text to stay a string
# (otherwise Excel auto-converts it to a boolean). That trick stamps the
# cell with a quote-prefix style, so strip formatting back to the default
# (unstyled) look by pasting formats in from a genuinely blank cell.
$blank = $ws.Range("H10")
$blank.Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
